$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: material text updated (perfil gslot -> gslot 2040), plus link
$ws.Range("A2").Value = "PERGIL GSLOT 2040 1 METRO"
$ws.Range("C2").Value = "https://www.makergal.es/product-page/perfil-gslot-2040-anodizado-natural"

# Rows 3-13, 16-17, 19-20: add ENLACE (link) values in column C
$ws.Range("C3").Value  = "https://amzn.to/3o6tNGl"
$ws.Range("C4").Value  = "https://amzn.to/3sH38U2"
$ws.Range("C5").Value  = "https://amzn.to/3iyzETC"
$ws.Range("C6").Value  = "https://amzn.to/398kDFf"
$ws.Range("C7").Value  = "https://amzn.to/3phdCHq"
$ws.Range("C8").Value  = "https://github.com/Srferrete/SFRSlider/tree/main/Stls"
$ws.Range("C9").Value  = "https://amzn.to/3634zCq"
$ws.Range("C10").Value = "https://www.makergal.es/product-page/polea-gt2-sincr%C3%B3nica-eje-3mm"
$ws.Range("C11").Value = "https://amzn.to/3sKe3MH"
$ws.Range("C12").Value = "https://www.makergal.es/product-page/ruedas-gslot"
$ws.Range("C13").Value = "https://www.makergal.es/product-page/tuerca-excentrica"
$ws.Range("C16").Value = "https://amzn.to/2KyGbBc"
$ws.Range("C17").Value = "https://amzn.to/3iENUdB"
$ws.Range("C19").Value = "https://amzn.to/3sLD7Tp"
$ws.Range("C20").Value = "https://amzn.to/3p4A5aP"

# Column C width widened to fit the URLs (~96 characters in real Excel ==
# stored width 95.28515625; this engine quantizes ColumnWidth to 1/6-character
# steps via a slightly different MDW model, so 94.5 is the closest input that
# lands on the nearest representable stored width)
$ws.Columns.Item(3).ColumnWidth = 94.5

# Selection moved to C21 (matches the author's last-edited cell)
$excel.Goto($ws.Range("C21"))
